$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (ISSCAAP 31)
$ws.Range("B6").Value = 1.113236270000001

# Row 10 (ISSCAAP 47)
$ws.Range("C10").Value = 1.35237506
$ws.Range("D10").Value = 6.341402055340052
$ws.Range("E10").Value = 39.51851520369351
$ws.Range("F10").Value = 54.14008274096643
$ws.Range("G10").Value = 45.85991725903357
$ws.Range("H10").Value = 54.14008274096643

# Row 11 (ISSCAAP 51)
$ws.Range("B11").Value = 4.579780390000003
$ws.Range("C11").Value = 4.324561374
$ws.Range("D11").Value = 25.99140387354132
$ws.Range("E11").Value = 44.04745699154409
$ws.Range("F11").Value = 29.9611391349146
$ws.Range("G11").Value = 70.03886086508541
$ws.Range("H11").Value = 29.9611391349146

# Row 12 (ISSCAAP 57)
$ws.Range("B12").Value = 5.56918511
$ws.Range("C12").Value = 5.47975862077466
$ws.Range("D12").Value = 14.93902192521486
$ws.Range("E12").Value = 67.50435457160448
$ws.Range("F12").Value = 17.55662350318067
$ws.Range("G12").Value = 82.44337649681934
$ws.Range("H12").Value = 17.55662350318067

# Row 15 (ISSCAAP 71)
$ws.Range("C15").Value = 9.088031675
$ws.Range("D15").Value = 22.93056935485755
$ws.Range("E15").Value = 26.57459035111507
$ws.Range("F15").Value = 50.49484029402739
$ws.Range("G15").Value = 49.50515970597262
$ws.Range("H15").Value = 50.49484029402739

# Row 16 (ISSCAAP 77)
$ws.Range("B16").Value = 1.544985800000001
$ws.Range("C16").Value = 1.324661300862069
$ws.Range("D16").Value = 43.27108330943184
$ws.Range("E16").Value = 38.29205176940727
$ws.Range("F16").Value = 18.43686492116091
$ws.Range("G16").Value = 81.56313507883911
$ws.Range("H16").Value = 18.43686492116091

# Row 17 (ISSCAAP 81)
$ws.Range("B17").Value = 0.3651133800000001

# Row 22 (Sharks) - fix double counts; B22 also changes number format to match C22 (3 decimals)
$ws.Range("B22").NumberFormat = "#,##0.000"
$ws.Range("B22").Value = 0.08399280000000002
$ws.Range("C22").Value = 0.05642497
$ws.Range("D22").Value = 48.68392486517937
$ws.Range("E22").Value = 37.54857113791996
$ws.Range("F22").Value = 13.76750399690066
$ws.Range("G22").Value = 86.23249600309933
$ws.Range("H22").Value = 13.76750399690066

# Row 24 (Global)
$ws.Range("B24").Value = 80.28049283000003
$ws.Range("C24").Value = 69.79157328945166
$ws.Range("D24").Value = 26.52088558783935
$ws.Range("E24").Value = 48.27323632967825
$ws.Range("F24").Value = 25.2058780824824
$ws.Range("G24").Value = 74.7941219175176
$ws.Range("H24").Value = 25.2058780824824
